# Actualización 11 de Mayo - Mañana
$wb = $excel.ActiveWorkbook

# ------------------------------------------------------------------
# Hoja "Estadisticos 1P"
# ------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Estadisticos 1P")
$ws1.Range("D2").Value = 12
$ws1.Range("F2").Value = 24
$ws1.Range("G2").Value = 66.67
$ws1.Range("H2").Value = 8.6

$ws1.Range("D3").Value = 7
$ws1.Range("F3").Value = 19
$ws1.Range("G3").Value = 73.08
$ws1.Range("H3").Value = 8.1

# ------------------------------------------------------------------
# Hoja "Estadisticos 2P"
# ------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Estadisticos 2P")
$ws2.Range("D2").Value = 17
$ws2.Range("E2").Value = 5
$ws2.Range("F2").Value = 19
$ws2.Range("G2").Value = 52.78
$ws2.Range("H2").Value = 9.2

$ws2.Range("D3").Value = 7
$ws2.Range("E3").Value = 0
$ws2.Range("F3").Value = 19
$ws2.Range("G3").Value = 73.08
$ws2.Range("H3").Value = 8.1

$ws2.Range("D5").Value = 15
$ws2.Range("E5").Value = 5
$ws2.Range("F5").Value = 7
$ws2.Range("G5").Value = 31.82
$ws2.Range("H5").Value = 7.6

# ------------------------------------------------------------------
# Hoja "Estadisticos Final"
# ------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("Estadisticos Final")
$ws3.Range("D2").Value = 12
$ws3.Range("F2").Value = 24
$ws3.Range("G2").Value = 66.67
$ws3.Range("H2").Value = 8.7

$ws3.Range("D3").Value = 7
$ws3.Range("F3").Value = 19
$ws3.Range("G3").Value = 73.08
$ws3.Range("H3").Value = 8.4

$ws3.Range("H5").Value = 7.4

# ------------------------------------------------------------------
# Hoja "Rescatables"
# ------------------------------------------------------------------
$ws4 = $wb.Worksheets.Item("Rescatables")

$data = @(
    @(2,  18330051920152, "CRISTOBAL",      "ROMERO",      "EDGAR ARMANDO",     "6AEM", 2),
    @(3,  18330051920159, "LOPEZ",          "ZAMUDIO",     "EZRA",              "6AEM", 2),
    @(4,  18330051920069, "MARIA",          "HERNANDEZ",   "AMALIO JARET",      "6AEM", 2),
    @(5,  18330051920170, "REYES",          "MARTINEZ",    "SALVADOR",          "6AEM", 2),
    @(6,  18330051920177, "SANDOVAL",       "GUZMAN",      "SAUL BRANDON",      "6AEM", 2),
    @(7,  18330051920176, "SANCHEZ",        "TRUJILLO",    "ERIK JAIR",         "6AEM", 2),
    @(8,  17330051920160, "RUIZ",           "LOPEZ",       "ALFONSO",           "6AEV", 2),
    @(9,  18330051920017, "DE LOS SANTOS",  "GONZALEZ",    "MARIA FERNANDA",    "6AEV", 2),
    @(10, 18330051920180, "VENEGAS",        "AMECA",       "ANGEL ISMAEL",      "6AEM", 1),
    @(11, 18330051920002, "ANASTACIO",      "HERNANDEZ",   "DIEGO APOLINAR",    "6AEV", 1),
    @(12, 18330051920003, "APARICIO",       "NAVARRO",     "PABLO",             "6AEV", 1),
    @(13, 18330051920004, "BAEZ",           "REYES",       "CRISTIAN MAURICIO", "6AEV", 1),
    @(14, 18330051920006, "BERNABE",        "NICIO",       "EMANUEL",           "6AEV", 1),
    @(15, 18330051920010, "COCOTLE",        "CUAQUEHUA",   "RIGOBERTO",         "6AEV", 1),
    @(16, 18330051920015, "COXCAHUA",       "TZITZIHUA",   "MARIA TERESA",      "6AEV", 1),
    @(17, 18330051920018, "DOMINGUEZ",      "TORRES",      "HUMBERTO",          "6AEV", 1),
    @(18, 18330051920019, "GALAN",          "CONSTANTINO", "JUAN CARLOS",       "6AEV", 1),
    @(19, 18330051920038, "JIMENEZ",        "DAMIAN",      "IRAIS",             "6AEV", 1),
    @(20, 18330051920024, "JIMENEZ",        "MOLOHUA",     "AZAEL",             "6AEV", 1),
    @(21, 18330051920027, "ROMAN",          "MARTINEZ",    "JULIO CESAR",       "6AEV", 1),
    @(22, 18330051920451, "TRUJILLO",       "ALVARADO",    "YAHIR ALEJANDRO",   "6AEV", 1),
    @(23, 18330051920036, "XOTLANIHUA",     "LORENZO",     "CARLOS MANUEL",     "6AEV", 1)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws4.Cells.Item($r, 1).Value = $row[1]
    $ws4.Cells.Item($r, 2).Value = $row[2]
    $ws4.Cells.Item($r, 3).Value = $row[3]
    $ws4.Cells.Item($r, 4).Value = $row[4]
    $ws4.Cells.Item($r, 6).Value = $row[5]
    $ws4.Cells.Item($r, 7).Value = $row[6]
}
